$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 126, shifting existing rows 126..231 down to 127..232.
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with the new data record.
$ws.Cells.Item(126, 1).Value = 5
$ws.Cells.Item(126, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(126, 3).Value = "Maule"
$ws.Cells.Item(126, 4).Value = 44669
$ws.Cells.Item(126, 4).NumberFormat = $ws.Cells.Item(127, 4).NumberFormat
$ws.Cells.Item(126, 5).Value = 7
$ws.Cells.Item(126, 6).Value = 100112009
$ws.Cells.Item(126, 7).Value = "Acelga"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 500
$ws.Cells.Item(126, 11).Value = 3500
$ws.Cells.Item(126, 12).Value = 3500
$ws.Cells.Item(126, 13).Value = 3500
$ws.Cells.Item(126, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(126, 15).Value = "Región del Maule"
$ws.Cells.Item(126, 16).Value = 875
$ws.Cells.Item(126, 17).Value = 4
$ws.Cells.Item(126, 18).Value = "Hortaliza"
